$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "최종점수" (K) and "MACRO_SCORE" (N) values for rows 2-6
$ws.Range("K2").Value = 55.1
$ws.Range("N2").Value = 52.28493729186943

$ws.Range("K3").Value = 51.1
$ws.Range("N3").Value = 52.28493729186943

$ws.Range("K4").Value = 46.9
$ws.Range("N4").Value = 52.28493729186943

$ws.Range("K5").Value = 45.7
$ws.Range("N5").Value = 52.28493729186943

$ws.Range("K6").Value = 37.1
$ws.Range("N6").Value = 52.28493729186943
